# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# The "Rules" sheet, cell B11 (row 11, col B) held the shared text "R40".
# It is re-typed as the literal text "1" (kept as text, not converted to
# the number 1) -- same cell, same row of the rule table, value swapped.
#
# A leading apostrophe is used so Excel treats the entry as text (quote
# prefix) instead of auto-coercing the numeric-looking "1" into a number;
# this keeps the cell a shared string, matching the edit that was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
